# Updated object list creation code. Updated code for the result container
# to add support for bus parameter export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix class name typo: PVsystems -> PVSystems
$ws.Range("A3").Value = "PVSystems"

# Add a new row for the Buses class / puVmagAngle property
$ws.Range("A4").Value = "Buses"
$ws.Range("B4").Value = "puVmagAngle"

# Copy the row-3 formatting down onto the new row so the new cells pick up
# the same table styling (font / vertical alignment) as the existing rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The new property cell was entered in white text (Consolas) - matches the
# authored workbook where this placeholder text isn't meant to stand out.
$f = $ws.Range("B4").Font
$f.Name = "Consolas"
$f.Family = 3
$f.ThemeColor = 2  # xlThemeColorLight1 -> theme="0" (white)

# Grow the Table1 ListObject so the new row is part of the table / autofilter.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A2:C4"))
